# Update presentation slide order.
#
# Original order:
#   1 Helping the Homeless with Qiosk
#   2 Why We Want to Help
#   3 What We Aim To Do
#   4 Our Team Structure
#   5 Our Implementation
#   6 Demo
#   7 Our Future Vision
#   8 References
#
# New order (the "Demo" slide moves earlier, right after "What We Aim To
# Do" and before "Our Team Structure" / "Our Implementation"):
#   1 Helping the Homeless with Qiosk
#   2 Why We Want to Help
#   3 What We Aim To Do
#   4 Demo
#   5 Our Team Structure
#   6 Our Implementation
#   7 Our Future Vision
#   8 References

$p = $ppt.ActivePresentation

# Locate the "Demo" slide by its title text rather than assuming a fixed
# index, then move it to slide position 4.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = ""
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $title = $shape.TextFrame.TextRange.Text
            break
        }
    }
    if ($title -eq "Demo") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $p.Slides.Item($targetIndex).MoveTo(4)
}
